$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.045.96'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '3.033.59'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.50%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.030.06'
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  +16.89%  '
$ws.Range('E11').Value = '  +4.94%  '
$ws.Range('E12').Value = '  +2.65%  '
$ws.Range('E13').Value = '  +3.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.30%  '
$ws.Range('D16').Value = '3.536.89'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('E17').Value = '  +3.96%  '
$ws.Range('D18').Value = '63.024.28'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('D19').Value = '3.032.18'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('E21').Value = '  +2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.698'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.24%  '
$ws.Range('E23').Value = '  +4.00%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.64%  '
$ws.Range('E26').Value = '  +8.68%  '
$ws.Range('E27').Value = '  +4.95%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.29'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +12.98%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.00%  '
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('E34').Value = '  +3.42%  '
$ws.Range('D35').Value = '0.0₃0863'
$ws.Range('E35').Value = '  +7.27%  '
$ws.Range('E36').Value = '  +2.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.46%  '
$ws.Range('E38').Value = '  +11.75%  '
$ws.Range('E39').Value = '  +9.46%  '
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.46'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  +1.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.311'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +17.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +16.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '395.16'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0361'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.44%  '
$ws.Range('D47').Value = '2.721.94'
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.69%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.23%  '
